$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B15: was stored as text "5" -> should become a real number 5
$ws.Range("B15").Value = 5

# New row 16 of annotation data
$ws.Range("A16").Value = "parisk"

# B16 must stay textual "2" (not auto-converted to a number), so force
# text formatting before assignment, then strip the style back off so no
# explicit style index ends up referenced on the cell.
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "2"
$ws.Range("B16").Style = "Normal"

$ws.Range("C16").Value = "无"
$ws.Range("D16").Value = "DIS"
$ws.Range("E16").Value = "WRI"
$ws.Range("F16").Value = "4cbdf296-0ef7-4a60-9d08-bf70fb941ab3"
$ws.Range("G16").Value = "SJTB5GZCb_annotated.xlsx"
$ws.Range("H16").Value = "The paper does not sufficiently discuss and compare the relevant neuroscience literature and related work."
